# Re-pull / push updated "dSF" (column F) values for several rows.
# Only column F values change; everything else on the sheet stays as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    10 = -6
    11 = 0
    13 = -1
    17 = -5
    25 = 6
    26 = -9
    31 = 3
    32 = -5
    37 = -4
    38 = -3
    43 = -5
    44 = -5
    45 = -6
    46 = -5
    53 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
